# edit.ps1 -- apply "finished approach 4! chapter 7" changes to the
# Index of Key Ideas, Insights and Contributions glossary lists.
#
# Strategy for each *new* glossary entry: locate an existing neighbouring
# paragraph (the "anchor"), Copy/Paste it so the new paragraph inherits the
# exact same pPr (numbering / Compact style) and rPr (bold term run)
# structure, then overwrite the term text and the trailing description
# text in place. This keeps the bold term / non-bold description split
# intact without hand-building run XML.

$d = $word.ActiveDocument

function Find-ParaIndex($doc, $prefix) {
    $paras = $doc.Paragraphs
    $n = $paras.Count
    for ($i = 1; $i -le $n; $i++) {
        $t = $paras.Item($i).Range.Text
        if ($t.StartsWith($prefix)) {
            return $i
        }
    }
    return -1
}

function Insert-GlossaryEntry($doc, $anchorTerm, $position, $newTerm, $newDesc) {
    # Find the existing paragraph whose bold term is exactly $anchorTerm.
    $idx = Find-ParaIndex $doc $anchorTerm
    if ($idx -eq -1) {
        throw "Insert-GlossaryEntry: anchor not found: $anchorTerm"
    }
    $anchor = $doc.Paragraphs.Item($idx)
    $anchorRange = $anchor.Range
    $anchorRange.Copy()

    if ($position -eq "Before") {
        $insertPos = $anchorRange.Start
        $newIdx = $idx
    } else {
        $insertPos = $anchorRange.End
        $newIdx = $idx + 1
    }
    $insertPoint = $doc.Range($insertPos, $insertPos)
    $insertPoint.Paste()

    # The pasted paragraph is an exact duplicate of the anchor (same
    # bold-term / space / description run split). Rewrite its text.
    $newPara = $doc.Paragraphs.Item($newIdx)
    $newRange = $newPara.Range
    $paraStart = $newRange.Start
    $anchorTermLen = $anchorTerm.Length

    $termRange = $doc.Range($paraStart, $paraStart + $anchorTermLen)
    $termRange.Text = $newTerm

    $newPara2 = $doc.Paragraphs.Item($newIdx)
    $newRange2 = $newPara2.Range
    $termEnd2 = $paraStart + $newTerm.Length
    # Paragraph.Range.Text includes the trailing paragraph mark (\r) --
    # exclude it from the replacement range.
    $restEnd = $newRange2.End - 1
    $restRange = $doc.Range($termEnd2, $restEnd)
    $restRange.Text = " " + $newDesc

    return $newIdx
}

# 1. "Life Interface" -> "Life Interface Design" (term-only rename).
$lifeIdx = Find-ParaIndex $d "Life Interface"
if ($lifeIdx -eq -1) {
    throw "Could not find 'Life Interface' glossary entry"
}
$lifePara = $d.Paragraphs.Item($lifeIdx)
$lifeRange = $lifePara.Range
$lifeTermRange = $d.Range($lifeRange.Start, $lifeRange.Start + "Life Interface".Length)
$lifeTermRange.Text = "Life Interface Design"

# 2. New entry "Pushing the Seams" before "Proxy Representations of Immobile Data".
Insert-GlossaryEntry $d "Proxy Representations of Immobile Data" "Before" `
    "Pushing the Seams" "- [ADD SECTIONREF]" | Out-Null

# 3. New entry "Surface Information Injustices" after "Shared Data Interaction".
Insert-GlossaryEntry $d "Shared Data Interaction" "After" `
    "Surface Information Injustices" "- [ADD SECTIONREF]" | Out-Null

# 4. New entry "Accessibility Tags (ARIA)" before "Activism".
Insert-GlossaryEntry $d "Activism" "Before" `
    "Accessibility Tags (ARIA)" "-" | Out-Null

# 5/6. New entries "Design, Adversarial" then "Design, Design After" before
#      "Design, Disrespectful" (each insert lands immediately before the
#      still-unmoved anchor, so order comes out Adversarial, Design After,
#      Disrespectful).
Insert-GlossaryEntry $d "Design, Disrespectful" "Before" `
    "Design, Adversarial" "-" | Out-Null
Insert-GlossaryEntry $d "Design, Disrespectful" "Before" `
    "Design, Design After" "-" | Out-Null

# 7. New entry "Device Tenancy" before "Digital Civics".
Insert-GlossaryEntry $d "Digital Civics" "Before" `
    "Device Tenancy" "-" | Out-Null

# 8. New entry "Empowerment in Use" before "Entities".
Insert-GlossaryEntry $d "Entities" "Before" `
    "Empowerment in Use" "-" | Out-Null

# 9. New entry "Information Landscape" before "Infrastructural Power, and its Four Levers".
Insert-GlossaryEntry $d "Infrastructural Power, and its Four Levers" "Before" `
    "Information Landscape" "-" | Out-Null

# 10. New entry "TrackerControl" before "Troubled Families".
Insert-GlossaryEntry $d "Troubled Families" "Before" `
    "TrackerControl" "- see Data Flow Auditing." | Out-Null

# 11. New entry "Web Extensions" after "Web Augmentation".
Insert-GlossaryEntry $d "Web Augmentation" "After" `
    "Web Extensions" "-" | Out-Null

Write-Host "Done."
